# Fan RPM measurement added
# This script reproduces, via Excel COM interop, the changes described by the
# target diff for texts.xlsx (TouchGFX Typography / Translation sheets).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Typography sheet
# ---------------------------------------------------------------------------
$typo = $wb.Worksheets.Item("Typography")

# Row 7 ("LCD_Default"): size 20 -> 28, and the fallback character is cleared.
$typo.Range("D7").Value = 28
$typo.Range("H7").ClearContents()

# Row 10: brand new "LCD_Small" typography entry (liquidcrystal.ttf, size 20,
# bpp 4, fallback "?" char, wildcard chars "012346789.-", widget wildcard
# "0123456789 :APM").
$typo.Range("B10").Value = "LCD_Small"
$typo.Range("C10").Value = "liquidcrystal.ttf"
$typo.Range("D10").Value = 20
$typo.Range("E10").Value = 4
$typo.Range("F10").Value = "?"
$typo.Range("G10").Value = "012346789.-"
$typo.Range("H10").Value = "0123456789 :APM"

# ---------------------------------------------------------------------------
# 2) Translation sheet
# ---------------------------------------------------------------------------
$tr = $wb.Worksheets.Item("Translation")

# Rows whose font switches from "LCD_Default" to the new "LCD_Small".
$tr.Range("C6").Value = "LCD_Small"
$tr.Range("C8").Value = "LCD_Small"
$tr.Range("C11").Value = "LCD_Small"

# Rows 16-41: re-write the whole translation table (font/alignment/value
# columns got reshuffled as labels were merged/removed and new ones added).
$rows = @(
    @{ Row = 16; B = "SingleUseId33"; C = "Small";       D = "Left";   F = "V" },
    @{ Row = 17; B = "SingleUseId34"; C = "LCD_Default"; D = "Right";  F = "<value>" },
    @{ Row = 18; B = "SingleUseId36"; C = "Tiny";        D = "Right";  F = "Power" },
    @{ Row = 19; B = "SingleUseId37"; C = "LCD_Large";   D = "Right";  F = "<value>" },
    @{ Row = 20; B = "SingleUseId38"; C = "Default";     D = "Left";   F = "W" },
    @{ Row = 21; B = "SingleUseId39"; C = "Small";       D = "Left";   F = "Charge" },
    @{ Row = 22; B = "SingleUseId40"; C = "Small";       D = "Left";   F = "Ah" },
    @{ Row = 23; B = "SingleUseId41"; C = "LCD_Small";   D = "Right";  F = "4234.234<value>" },
    @{ Row = 24; B = "SingleUseId42"; C = "Tiny";        D = "Left";   F = "Capacity" },
    @{ Row = 25; B = "SingleUseId45"; C = "Tiny";        D = "Left";   F = "Elapsed time" },
    @{ Row = 26; B = "SingleUseId46"; C = "LCD_Small";   D = "Right";  F = "<>" },
    @{ Row = 27; B = "SingleUseId47"; C = "Small";       D = "Right";  F = "Active Load 8A max" },
    @{ Row = 28; B = "SingleUseId48"; C = "Small";       D = "Left";   F = "%" },
    @{ Row = 29; B = "off";           C = "LCD_Small";   D = "Left";   F = "off" },
    @{ Row = 30; B = "SingleUseId49"; C = "Small";       D = "Center"; F = "Reset" },
    @{ Row = 31; B = "SingleUseId50"; C = "Small";       D = "Left";   F = "Load Limits" },
    @{ Row = 32; B = "SingleUseId51"; C = "Small";       D = "Left";   F = "A" },
    @{ Row = 33; B = "SingleUseId52"; C = "LCD_Default"; D = "Right";  F = "8.250<value>" },
    @{ Row = 34; B = "SingleUseId54"; C = "Small";       D = "Left";   F = "V" },
    @{ Row = 35; B = "SingleUseId57"; C = "Small";       D = "Center"; F = "Start" },
    @{ Row = 36; B = "SingleUseId64"; C = "Small";       D = "Center"; F = "Yes" },
    @{ Row = 37; B = "SingleUseId65"; C = "Small";       D = "Center"; F = "Cancel" },
    @{ Row = 38; B = "SingleUseId66"; C = "Default";     D = "Left";   F = "Reset timer?" },
    @{ Row = 39; B = "SingleUseId67"; C = "Tiny";        D = "Right";  F = "Enable" },
    @{ Row = 40; B = "voltageCurrent";C = "Tiny";        D = "Left";   F = "Voltage / Current" },
    @{ Row = 41; B = "SingleUseId68"; C = "LCD_Default"; D = "Right";  F = "23.4" }
)

foreach ($item in $rows) {
    $r = $item.Row
    $tr.Cells.Item($r, 2).Value = $item.B   # column B
    $tr.Cells.Item($r, 3).Value = $item.C   # column C
    $tr.Cells.Item($r, 4).Value = $item.D   # column D
    $tr.Cells.Item($r, 5).Value = "LTR"     # column E (unchanged, re-asserted)
    $tr.Cells.Item($r, 6).Value = $item.F   # column F
}

# Rows 42-44 previously held translation entries that have now been merged
# away / removed entirely - clear them out completely.
$tr.Range("B42:F42").ClearContents()
$tr.Range("B43:F43").ClearContents()
$tr.Range("B44:F44").ClearContents()
